$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1685.7142
$ws.Range("I40").Value = 1600
$ws.Range("J40").Value = 1800
$ws.Range("K40").Value = 1600
$ws.Range("L40").Value = 1800
$ws.Range("M40").Value = -1425
$ws.Range("N40").Value = -2150
$ws.Range("H62").Value = 4637.5
$ws.Range("I62").Value = 3871.4285
$ws.Range("K62").Value = 3871.4285
$ws.Range("M62").Value = -3247.4285
$ws.Range("H65").Value = 4637.5
$ws.Range("I65").Value = 3871.4285
$ws.Range("K65").Value = 19357.1425
$ws.Range("M65").Value = -16237.1425
$ws.Range("H76").Value = 4578.222
$ws.Range("I76").Value = 5000
$ws.Range("K76").Value = 5000
$ws.Range("M76").Value = -4685
$ws.Range("H79").Value = 4578.222
$ws.Range("I79").Value = 5000
$ws.Range("K79").Value = 5000
$ws.Range("M79").Value = -3908
$ws.Range("H92").Value = 3555.3635
$ws.Range("I92").Value = 3767.6667
$ws.Range("J92").Value = 2600
$ws.Range("K92").Value = 3767.6667
$ws.Range("L92").Value = 2600
$ws.Range("M92").Value = -2519.6667
$ws.Range("N92").Value = -5096
$ws.Range("H129").Value = 831.69385
$ws.Range("J129").Value = 1019
$ws.Range("L129").Value = 3057
$ws.Range("N129").Value = -13057
$ws.Range("H132").Value = 4634435
$ws.Range("I132").Value = 5955575
$ws.Range("J132").Value = 10444.5625
$ws.Range("K132").Value = 17866725
$ws.Range("L132").Value = 31333.6875
$ws.Range("M132").Value = -17864195
$ws.Range("N132").Value = -36393.6875
$ws.Range("H137").Value = 1026.8276
$ws.Range("I137").Value = 844.6875
$ws.Range("J137").Value = 1251
$ws.Range("K137").Value = 2534.0625
$ws.Range("L137").Value = 3753
$ws.Range("M137").Value = 15.9375
$ws.Range("N137").Value = -8853
$ws.Range("H138").Value = 1214.45
$ws.Range("I138").Value = 503.7619
$ws.Range("J138").Value = 1729.0862
$ws.Range("K138").Value = 1511.2857
$ws.Range("L138").Value = 5187.2586
$ws.Range("M138").Value = 3628.7143
$ws.Range("N138").Value = -15467.2586
$ws.Range("H141").Value = 623.55554
$ws.Range("I141").Value = 520.9216
$ws.Range("J141").Value = 2368.3333
$ws.Range("K141").Value = 1562.7648
$ws.Range("L141").Value = 7104.999899999999
$ws.Range("M141").Value = 3617.2352
$ws.Range("N141").Value = -17464.9999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 595.1818
$ws.Range("I97").Value = 471.8889
$ws.Range("J97").Value = 1150
$ws.Range("K97").Value = 471.8889
$ws.Range("L97").Value = 1150
$ws.Range("M97").Value = 24.11110000000002
$ws.Range("N97").Value = -2142
$ws.Range("H132").Value = 1649.619
$ws.Range("I132").Value = 1533.1522
$ws.Range("J132").Value = 1964.7646
$ws.Range("K132").Value = 4599.4566
$ws.Range("L132").Value = 5894.293799999999
$ws.Range("M132").Value = -2069.4566
$ws.Range("N132").Value = -10954.2938

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 41669164
$ws.Range("I94").Value = 50002396
$ws.Range("K94").Value = 50002396
$ws.Range("M94").Value = -50001945
$ws.Range("H99").Value = 40001052
$ws.Range("I99").Value = 58824296
$ws.Range("J99").Value = 1649.875
$ws.Range("K99").Value = 58824296
$ws.Range("L99").Value = 1649.875
$ws.Range("M99").Value = -58822798
$ws.Range("N99").Value = -4645.875
$ws.Range("H105").Value = 62501680
$ws.Range("I105").Value = 76924730
$ws.Range("K105").Value = 76924730
$ws.Range("M105").Value = -76922983

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 90910616
$ws.Range("I16").Value = 90910616
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 90910616
$ws.Range("L16").Value = 0
$ws.Range("M16").ClearContents()
$ws.Range("N16").Value = -90910329
$ws.Range("H31").Value = 2012.75
$ws.Range("I31").Value = 2035.0435
$ws.Range("J31").Value = 1500
$ws.Range("K31").Value = 2035.0435
$ws.Range("L31").Value = 1500
$ws.Range("M31").Value = -1740.0435
$ws.Range("N31").Value = -2090
$ws.Range("H34").Value = 2012.75
$ws.Range("I34").Value = 2035.0435
$ws.Range("J34").Value = 1500
$ws.Range("K34").Value = 2035.0435
$ws.Range("L34").Value = 1500
$ws.Range("M34").Value = -1833.0435
$ws.Range("N34").Value = -1904
$ws.Range("H113").Value = 90910616
$ws.Range("I113").Value = 90910616
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 90910616
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -90908446
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").ClearContents()
$ws.Range("N135").Value = 0

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1963.3684
$ws.Range("I5").Value = 2228.8572
$ws.Range("K5").Value = 6686.571599999999
$ws.Range("M5").Value = -6574.571599999999
$ws.Range("H14").Value = 196.35294
$ws.Range("I14").Value = 196.35294
$ws.Range("K14").Value = 589.05882
$ws.Range("M14").Value = -416.05882
$ws.Range("H23").Value = 999.75
$ws.Range("I23").Value = 1400
$ws.Range("J23").Value = 759.6
$ws.Range("K23").Value = 4200
$ws.Range("L23").Value = 2278.8
$ws.Range("M23").Value = -3965
$ws.Range("N23").Value = -2748.8
$ws.Range("H32").Value = 1925.1666
$ws.Range("I32").Value = 751
$ws.Range("J32").Value = 2160
$ws.Range("K32").Value = 2253
$ws.Range("L32").Value = 6480
$ws.Range("M32").Value = -1970
$ws.Range("N32").Value = -7046
$ws.Range("H40").Value = 231.2
$ws.Range("I40").Value = 108.625
$ws.Range("K40").Value = 434.5
$ws.Range("M40").Value = -365.5
$ws.Range("H86").Value = 690
$ws.Range("I86").Value = 690
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2070
$ws.Range("L86").ClearContents()
$ws.Range("M86").Value = -884
$ws.Range("N86").Value = 0
$ws.Range("H89").Value = 690
$ws.Range("I89").Value = 690
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 6210
$ws.Range("L89").ClearContents()
$ws.Range("M89").Value = -282
$ws.Range("N89").Value = 0
$ws.Range("H113").Value = 699.03125
$ws.Range("J113").Value = 705.4516
$ws.Range("L113").Value = 2116.3548
$ws.Range("N113").Value = -6456.3548
$ws.Range("H122").Value = 823
$ws.Range("I122").Value = 640.75
$ws.Range("J122").Value = 889.2727
$ws.Range("K122").Value = 5766.75
$ws.Range("L122").Value = 8003.454299999999
$ws.Range("M122").Value = -3316.75
$ws.Range("N122").Value = -12903.4543
$ws.Range("H131").Value = 19609110
$ws.Range("I131").Value = 100000440
$ws.Range("J131").Value = 1469.7317
$ws.Range("K131").Value = 300001320
$ws.Range("L131").Value = 4409.1951
$ws.Range("M131").Value = -299996280
$ws.Range("N131").Value = -14489.1951
$ws.Range("H132").Value = 1829.8572
$ws.Range("J132").Value = 1535
$ws.Range("L132").Value = 13815
$ws.Range("N132").Value = -18875
$ws.Range("H135").Value = 1963.3684
$ws.Range("I135").Value = 2228.8572
$ws.Range("K135").Value = 20059.7148
$ws.Range("M135").Value = -17524.7148
$ws.Range("H139").Value = 1424.8
$ws.Range("I139").Value = 1424.8
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 4274.4
$ws.Range("L139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = 865.6000000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 0
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H130").Value = 37340
$ws.Range("J130").Value = 37340
$ws.Range("L130").Value = 37340
$ws.Range("N130").Value = -47380

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 39800
$ws.Range("J36").Value = 39800
$ws.Range("L36").Value = 39800
$ws.Range("N36").Value = -40924
$ws.Range("H93").Value = 761.6
$ws.Range("I93").Value = 761.6
$ws.Range("K93").Value = 761.6
$ws.Range("M93").Value = 486.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4497.5
$ws.Range("J96").Value = 2992.5
$ws.Range("L96").Value = 2992.5
$ws.Range("N96").Value = -5738.5
$ws.Range("H100").Value = 675.4286
$ws.Range("J100").Value = 609.5
$ws.Range("L100").Value = 1219
$ws.Range("N100").Value = -2301
$ws.Range("H107").Value = 401.3
$ws.Range("I107").Value = 375.05884
$ws.Range("J107").Value = 550
$ws.Range("K107").Value = 1125.17652
$ws.Range("L107").Value = 1650
$ws.Range("M107").Value = 794.82348
$ws.Range("N107").Value = -5490
$ws.Range("H126").Value = 40001256
$ws.Range("I126").Value = 50001000
$ws.Range("J126").Value = 2280
$ws.Range("K126").Value = 150003000
$ws.Range("L126").Value = 6840
$ws.Range("M126").Value = -150000530
$ws.Range("N126").Value = -11780
